# Faraday Lab data tables.xlsx - apply commit "added field to excel"
#
# Summary of the edit:
#  - Sheet3 ("changeCurrent") gains two new columns F ("B (mT)") and G ("dB"),
#    interleaved literal values / AVERAGE() formulas, a custom width on column B,
#    and becomes the active/selected sheet.
#  - Sheet3 headers B1/D1 get units appended ("B*L" -> "B*L (mT*cm)", "V" -> "V (V)").
#  - Sheet1 ("ChangeThetaRaw") stops being the tab-selected sheet.
#  - Sheet2 ("changeThetaCalculations") selection moves to H2:I6.
#  - Sheet3 selection moves to H16.
#  - Workbook active tab becomes sheet3 (index 2, 0-based).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Sheet3 ("changeCurrent"): rename headers with units ---------------------
$ws3.Range("B1").Value = "B*L (mT*cm)"
$ws3.Range("D1").Value = "V (V)"

# --- Sheet3: new header cells for the two added columns ----------------------
$ws3.Range("F1").Value = "B (mT)"
$ws3.Range("G1").Value = "dB"

# --- Sheet3: column B custom width (~13.16 chars) -----------------------------
$ws3.Columns.Item(2).ColumnWidth = 12.33

# --- Sheet3: fill column F ("B (mT)") - literals on even measured rows, -------
# AVERAGE() of neighbours on the interpolated rows ----------------------------
$ws3.Range("F2").Value = -31.8
$ws3.Range("F3").Formula = "=AVERAGE(F2,F4)"
$ws3.Range("F4").Value = -21.2
$ws3.Range("F5").Formula = "=AVERAGE(F4,F6)"
$ws3.Range("F6").Value = -10.6
$ws3.Range("F7").Formula = "=AVERAGE(F6,F8)"
$ws3.Range("F8").Value = 0
$ws3.Range("F9").Formula = "=AVERAGE(F8,F10)"
$ws3.Range("F10").Value = 10.6
$ws3.Range("F11").Formula = "=AVERAGE(F10,F12)"
$ws3.Range("F12").Value = 21.2
$ws3.Range("F13").Formula = "=AVERAGE(F12,F14)"
$ws3.Range("F14").Value = 31.8

# --- Sheet3: fill column G ("dB") ---------------------------------------------
$ws3.Range("G2").Value = -1
$ws3.Range("G3").Value = -0.9
$ws3.Range("G4").Value = -0.7
$ws3.Range("G5").Formula = "=AVERAGE(G4,G6)"
$ws3.Range("G6").Value = -0.3
$ws3.Range("G7").Formula = "=AVERAGE(G6,G8)"
$ws3.Range("G8").Value = 0
$ws3.Range("G9").Formula = "=AVERAGE(G8,G10)"
$ws3.Range("G10").Value = 0.3
$ws3.Range("G11").Formula = "=AVERAGE(G10,G12)"
$ws3.Range("G12").Value = 0.7
$ws3.Range("G13").Value = 0.9
$ws3.Range("G14").Value = 1

# --- Sheet2 ("changeThetaCalculations"): update selection ---------------------
$ws2.Activate()
[void]$ws2.Range("H2:I6").Select()

# --- Sheet3 ("changeCurrent"): becomes the active/selected sheet + selection --
$ws3.Activate()
[void]$ws3.Range("H16").Select()
